# Fruta / hortaliza, semanal
# The weekly refresh re-shuffles the "Fecha" (D), "Volumen" (J),
# "Precio minimo" (K), "Precio maximo" (L), "Precio promedio ponderado" (M)
# and "Precio $/Kg" (P) figures across the existing data rows (2-12, row 7
# untouched). Capture the original values first so the row-to-row copy
# doesn't clobber a source row before it has been read.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# columns that move together as one record
$cols = @(4, 10, 11, 12, 13, 16)   # D, J, K, L, M, P

# target row -> source row (source row's ORIGINAL values land in target row)
$rowMap = @{
    2  = 6
    3  = 9
    4  = 11
    5  = 12
    6  = 3
    8  = 4
    9  = 10
    10 = 5
    11 = 2
    12 = 8
}

# snapshot the original values for every row referenced above
$orig = @{}
$rowsNeeded = @(2, 3, 4, 5, 6, 8, 9, 10, 11, 12)
foreach ($r in $rowsNeeded) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $orig[$r] = $rowVals
}

# write the shuffled values back using the captured snapshot
foreach ($targetRow in $rowMap.Keys) {
    $sourceRow = $rowMap[$targetRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($targetRow, $c).Value = $orig[$sourceRow][$c]
    }
}
